$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 758.05884
$ws.Range("I43").Value = 582.1111
$ws.Range("J43").Value = 956
$ws.Range("K43").Value = 582.1111
$ws.Range("L43").Value = 956
$ws.Range("M43").Value = -513.1111
$ws.Range("N43").Value = -1094

$ws.Range("H64").Value = 3607.8875
$ws.Range("I64").Value = 3359.9778
$ws.Range("J64").Value = 3926.6287
$ws.Range("K64").Value = 3359.9778
$ws.Range("L64").Value = 3926.6287
$ws.Range("M64").Value = -3111.9778
$ws.Range("N64").Value = -4422.6287

$ws.Range("H67").Value = 3607.8875
$ws.Range("I67").Value = 3359.9778
$ws.Range("J67").Value = 3926.6287
$ws.Range("K67").Value = 3359.9778
$ws.Range("L67").Value = 3926.6287
$ws.Range("M67").Value = -2501.9778
$ws.Range("N67").Value = -5642.6287

$ws.Range("H74").Value = 3417.1936
$ws.Range("I74").Value = 3100
$ws.Range("J74").Value = 3478.1924
$ws.Range("K74").Value = 3100
$ws.Range("L74").Value = 3478.1924
$ws.Range("M74").Value = -2164
$ws.Range("N74").Value = -5350.1924

$ws.Range("H76").Value = 2977.9285
$ws.Range("I76").Value = 2555.3333
$ws.Range("J76").Value = 3178.1052
$ws.Range("K76").Value = 2555.3333
$ws.Range("L76").Value = 3178.1052
$ws.Range("M76").Value = -2240.3333
$ws.Range("N76").Value = -3808.1052

$ws.Range("H77").Value = 3417.1936
$ws.Range("I77").Value = 3100
$ws.Range("J77").Value = 3478.1924
$ws.Range("K77").Value = 15500
$ws.Range("L77").Value = 17390.962
$ws.Range("M77").Value = -10820
$ws.Range("N77").Value = -26750.962

$ws.Range("H79").Value = 2977.9285
$ws.Range("I79").Value = 2555.3333
$ws.Range("J79").Value = 3178.1052
$ws.Range("K79").Value = 2555.3333
$ws.Range("L79").Value = 3178.1052
$ws.Range("M79").Value = -1463.3333
$ws.Range("N79").Value = -5362.1052

$ws.Range("H129").Value = 778.87756
$ws.Range("I129").Value = 350
$ws.Range("J129").Value = 817
$ws.Range("K129").Value = 1050
$ws.Range("L129").Value = 2451
$ws.Range("M129").Value = 3950
$ws.Range("N129").Value = -12451

$ws.Range("H132").Value = 2578.1482
$ws.Range("I132").Value = 2118.8723
$ws.Range("K132").Value = 6356.6169
$ws.Range("M132").Value = -3826.6169

$ws.Range("H135").Value = 418.69232
$ws.Range("I135").Value = 326.6316
$ws.Range("J135").Value = 668.5714
$ws.Range("K135").Value = 2939.6844
$ws.Range("L135").Value = 6017.1426
$ws.Range("M135").Value = -404.6844000000001
$ws.Range("N135").Value = -11087.1426

$ws.Range("H137").Value = 12501610
$ws.Range("I137").Value = 20000820
$ws.Range("K137").Value = 60002460
$ws.Range("M137").Value = -59999910

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1504.6
$ws.Range("J88").Value = 1621.4286
$ws.Range("L88").Value = 1621.4286
$ws.Range("N88").Value = -2433.4286

$ws.Range("H91").Value = 1504.6
$ws.Range("J91").Value = 1621.4286
$ws.Range("L91").Value = 1621.4286
$ws.Range("N91").Value = -4429.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2969.8655
$ws.Range("I105").Value = 1800
$ws.Range("J105").Value = 3041.4897
$ws.Range("K105").Value = 1800
$ws.Range("L105").Value = 3041.4897
$ws.Range("M105").Value = -53
$ws.Range("N105").Value = -6535.4897

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2935
$ws.Range("I31").Value = 1895.85
$ws.Range("J31").Value = 5904
$ws.Range("K31").Value = 1895.85
$ws.Range("L31").Value = 5904
$ws.Range("M31").Value = -1600.85
$ws.Range("N31").Value = -6494

$ws.Range("H34").Value = 2935
$ws.Range("I34").Value = 1895.85
$ws.Range("J34").Value = 5904
$ws.Range("K34").Value = 1895.85
$ws.Range("L34").Value = 5904
$ws.Range("M34").Value = -1693.85
$ws.Range("N34").Value = -6308

$ws.Range("H58").Value = 2411.2693
$ws.Range("I58").Value = 1212.4445
$ws.Range("J58").Value = 5108.625
$ws.Range("K58").Value = 1212.4445
$ws.Range("L58").Value = 5108.625
$ws.Range("M58").Value = -1009.4445
$ws.Range("N58").Value = -5514.625

$ws.Range("H136").Value = 2411.2693
$ws.Range("I136").Value = 1212.4445
$ws.Range("J136").Value = 5108.625
$ws.Range("K136").Value = 3637.3335
$ws.Range("L136").Value = 15325.875
$ws.Range("M136").Value = -1087.3335
$ws.Range("N136").Value = -20425.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 4812.75
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 5214.5713
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 15643.7139
$ws.Range("M17").Value = -5831
$ws.Range("N17").Value = -15981.7139

$ws.Range("H93").Value = 4984.95
$ws.Range("J93").Value = 4984.95
$ws.Range("L93").Value = 14954.85
$ws.Range("N93").Value = -18698.85

$ws.Range("H108").Value = 2804.6155
$ws.Range("I108").Value = 1814.2858
$ws.Range("J108").Value = 3960
$ws.Range("K108").Value = 5442.857400000001
$ws.Range("L108").Value = 11880
$ws.Range("M108").Value = -2562.857400000001
$ws.Range("N108").Value = -17640

$ws.Range("H120").Value = 15850
$ws.Range("J120").Value = 19350
$ws.Range("L120").Value = 58050
$ws.Range("N120").Value = -67726

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H70").Value = 6963.591
$ws.Range("I70").Value = 7494.4443
$ws.Range("J70").Value = 4574.75
$ws.Range("K70").Value = 7494.4443
$ws.Range("L70").Value = 4574.75
$ws.Range("M70").Value = -7224.4443
$ws.Range("N70").Value = -5114.75

$ws.Range("H73").Value = 6963.591
$ws.Range("I73").Value = 7494.4443
$ws.Range("J73").Value = 4574.75
$ws.Range("K73").Value = 7494.4443
$ws.Range("L73").Value = 4574.75
$ws.Range("M73").Value = -6558.4443
$ws.Range("N73").Value = -6446.75

$ws.Range("H80").Value = 3000.75
$ws.Range("I80").Value = 2818.182
$ws.Range("J80").Value = 3155.2307
$ws.Range("K80").Value = 2818.182
$ws.Range("L80").Value = 3155.2307
$ws.Range("M80").Value = -1820.182
$ws.Range("N80").Value = -5151.2307

$ws.Range("H83").Value = 3000.75
$ws.Range("I83").Value = 2818.182
$ws.Range("J83").Value = 3155.2307
$ws.Range("K83").Value = 14090.91
$ws.Range("L83").Value = 15776.1535
$ws.Range("M83").Value = -9098.91
$ws.Range("N83").Value = -25760.1535

$ws.Range("H102").Value = 1763.75
$ws.Range("I102").Value = 1530
$ws.Range("J102").Value = 3400
$ws.Range("K102").Value = 1530
$ws.Range("L102").Value = 3400
$ws.Range("M102").Value = 92
$ws.Range("N102").Value = -6644

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 670.61536
$ws.Range("I22").Value = 562.5
$ws.Range("J22").Value = 718.6667
$ws.Range("K22").Value = 562.5
$ws.Range("L22").Value = 718.6667
$ws.Range("M22").Value = -267.5
$ws.Range("N22").Value = -1308.6667

$ws.Range("H27").Value = 670.61536
$ws.Range("I27").Value = 562.5
$ws.Range("J27").Value = 718.6667
$ws.Range("K27").Value = 562.5
$ws.Range("L27").Value = 718.6667
$ws.Range("M27").Value = -455.5
$ws.Range("N27").Value = -932.6667

$ws.Range("H44").Value = 4800
$ws.Range("J44").Value = 4800
$ws.Range("L44").Value = 4800
$ws.Range("N44").Value = -5712

$ws.Range("H122").Value = 3319.8438
$ws.Range("I122").Value = 2974.2
$ws.Range("J122").Value = 4554.2856
$ws.Range("K122").Value = 8922.599999999999
$ws.Range("L122").Value = 13662.8568
$ws.Range("M122").Value = -6472.599999999999
$ws.Range("N122").Value = -18562.8568

$ws.Range("H136").Value = 1834.8148
$ws.Range("I136").Value = 1076.5
$ws.Range("J136").Value = 4001.4285
$ws.Range("K136").Value = 3229.5
$ws.Range("L136").Value = 12004.2855
$ws.Range("M136").Value = -679.5
$ws.Range("N136").Value = -17104.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3670.6667
$ws.Range("I132").Value = 3585.8948
$ws.Range("J132").Value = 3765.4119
$ws.Range("K132").Value = 10757.6844
$ws.Range("L132").Value = 11296.2357
$ws.Range("M132").Value = -8227.6844
$ws.Range("N132").Value = -16356.2357

Write-Output "Edit complete"
